$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1893.25
$ws.Range("I88").Value = 699
$ws.Range("J88").Value = 2291.3333
$ws.Range("K88").Value = 699
$ws.Range("L88").Value = 2291.3333
$ws.Range("M88").Value = -293
$ws.Range("N88").Value = -3103.3333
$ws.Range("H91").Value = 1893.25
$ws.Range("I91").Value = 699
$ws.Range("J91").Value = 2291.3333
$ws.Range("K91").Value = 699
$ws.Range("L91").Value = 2291.3333
$ws.Range("M91").Value = 705
$ws.Range("N91").Value = -5099.3333
$ws.Range("H106").Value = 5187.5
$ws.Range("I106").Value = 4583.3335
$ws.Range("K106").Value = 4583.3335
$ws.Range("M106").Value = -3952.3335
$ws.Range("H107").Value = 564.5
$ws.Range("I107").Value = 468.75
$ws.Range("J107").Value = 660.25
$ws.Range("K107").Value = 468.75
$ws.Range("L107").Value = 660.25
$ws.Range("M107").Value = 1451.25
$ws.Range("N107").Value = -4500.25
$ws.Range("H111").Value = 39666.332
$ws.Range("I111").Value = 9000
$ws.Range("J111").Value = 54999.5
$ws.Range("K111").Value = 27000
$ws.Range("L111").Value = 164998.5
$ws.Range("M111").Value = -23933
$ws.Range("N111").Value = -171132.5
$ws.Range("H125").Value = 2915
$ws.Range("I125").Value = 3649.75
$ws.Range("K125").Value = 32847.75
$ws.Range("M125").Value = -30387.75
$ws.Range("H127").Value = 740.5714
$ws.Range("I127").Value = 740.5714
$ws.Range("K127").Value = 2221.7142
$ws.Range("M127").Value = 2738.2858
$ws.Range("H138").Value = 3147.7273
$ws.Range("I138").Value = 1571.12
$ws.Range("J138").Value = 4461.567
$ws.Range("K138").Value = 4713.36
$ws.Range("L138").Value = 13384.701
$ws.Range("M138").Value = 426.6400000000003
$ws.Range("N138").Value = -23664.701

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 2401.5
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 4800
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 4800
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -5146
$ws.Range("H32").Value = 16695.113
$ws.Range("I32").Value = 21472.516
$ws.Range("J32").Value = 2362.9092
$ws.Range("K32").Value = 21472.516
$ws.Range("L32").Value = 2362.9092
$ws.Range("M32").Value = -21185.516
$ws.Range("N32").Value = -2936.9092
$ws.Range("H61").Value = 4130
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 3757.1428
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 3757.1428
$ws.Range("M61").Value = -4788
$ws.Range("N61").Value = -4181.1428
$ws.Range("H74").Value = 1414.3529
$ws.Range("I74").Value = 1945.5
$ws.Range("J74").Value = 942.2222
$ws.Range("K74").Value = 1945.5
$ws.Range("L74").Value = 942.2222
$ws.Range("M74").Value = -1071.5
$ws.Range("N74").Value = -2690.2222
$ws.Range("H77").Value = 1414.3529
$ws.Range("I77").Value = 1945.5
$ws.Range("J77").Value = 942.2222
$ws.Range("K77").Value = 9727.5
$ws.Range("L77").Value = 4711.111
$ws.Range("M77").Value = -5359.5
$ws.Range("N77").Value = -13447.111
$ws.Range("H122").Value = 6151.75
$ws.Range("I122").Value = 7369
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 22107
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -19657
$ws.Range("N122").Value = -12400
$ws.Range("H136").Value = 4130
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 3757.1428
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 11271.4284
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -16371.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 52504.1
$ws.Range("I86").Value = 2653.5334
$ws.Range("J86").Value = 202055.8
$ws.Range("K86").Value = 2653.5334
$ws.Range("L86").Value = 202055.8
$ws.Range("M86").Value = -1530.5334
$ws.Range("N86").Value = -204301.8
$ws.Range("H89").Value = 52504.1
$ws.Range("I89").Value = 2653.5334
$ws.Range("J89").Value = 202055.8
$ws.Range("K89").Value = 13267.667
$ws.Range("L89").Value = 1010279
$ws.Range("M89").Value = -7651.666999999999
$ws.Range("N89").Value = -1021511
$ws.Range("H107").Value = 22454.6
$ws.Range("I107").Value = 29564.889
$ws.Range("J107").Value = 4171
$ws.Range("K107").Value = 29564.889
$ws.Range("L107").Value = 4171
$ws.Range("M107").Value = -27644.889
$ws.Range("N107").Value = -8011
$ws.Range("H134").Value = 2277.6216
$ws.Range("I134").Value = 2008.5
$ws.Range("K134").Value = 6025.5
$ws.Range("M134").Value = -3490.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1990.6666
$ws.Range("I5").Value = 20004
$ws.Range("J5").Value = 1090
$ws.Range("K5").Value = 60012
$ws.Range("L5").Value = 3270
$ws.Range("M5").Value = -59900
$ws.Range("N5").Value = -3494
$ws.Range("H56").Value = 4605.58
$ws.Range("I56").Value = 4605.58
$ws.Range("K56").Value = 4605.58
$ws.Range("M56").Value = -4075.58
$ws.Range("H135").Value = 1990.6666
$ws.Range("I135").Value = 20004
$ws.Range("J135").Value = 1090
$ws.Range("K135").Value = 180036
$ws.Range("L135").Value = 9810
$ws.Range("M135").Value = -177501
$ws.Range("N135").Value = -14880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7573.6924
$ws.Range("I70").Value = 7975.8
$ws.Range("J70").Value = 6233.3335
$ws.Range("K70").Value = 7975.8
$ws.Range("L70").Value = 6233.3335
$ws.Range("M70").Value = -7705.8
$ws.Range("N70").Value = -6773.3335
$ws.Range("H73").Value = 7573.6924
$ws.Range("I73").Value = 7975.8
$ws.Range("J73").Value = 6233.3335
$ws.Range("K73").Value = 7975.8
$ws.Range("L73").Value = 6233.3335
$ws.Range("M73").Value = -7039.8
$ws.Range("N73").Value = -8105.3335
$ws.Range("H80").Value = 3829.3572
$ws.Range("I80").Value = 3767.2222
$ws.Range("J80").Value = 3941.2
$ws.Range("K80").Value = 3767.2222
$ws.Range("L80").Value = 3941.2
$ws.Range("M80").Value = -2769.2222
$ws.Range("N80").Value = -5937.2
$ws.Range("H83").Value = 3829.3572
$ws.Range("I83").Value = 3767.2222
$ws.Range("J83").Value = 3941.2
$ws.Range("K83").Value = 18836.111
$ws.Range("L83").Value = 19706
$ws.Range("M83").Value = -13844.111
$ws.Range("N83").Value = -29690
$ws.Range("H97").Value = 40727.46
$ws.Range("I97").Value = 57117.445
$ws.Range("J97").Value = 3850
$ws.Range("K97").Value = 57117.445
$ws.Range("L97").Value = 3850
$ws.Range("M97").Value = -56621.445
$ws.Range("N97").Value = -4842
$ws.Range("H102").Value = 2206.1365
$ws.Range("I102").Value = 1789.7142
$ws.Range("J102").Value = 2934.875
$ws.Range("K102").Value = 1789.7142
$ws.Range("L102").Value = 2934.875
$ws.Range("M102").Value = -167.7141999999999
$ws.Range("N102").Value = -6178.875
$ws.Range("H122").Value = 4396.364
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4436
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 13308
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -18208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4773.25
$ws.Range("I7").Value = 5307.7
$ws.Range("J7").Value = 3882.5
$ws.Range("K7").Value = 5307.7
$ws.Range("L7").Value = 3882.5
$ws.Range("M7").Value = -5195.7
$ws.Range("N7").Value = -4106.5
$ws.Range("H122").Value = 64288196
$ws.Range("I122").Value = 62502144
$ws.Range("J122").Value = 66669600
$ws.Range("K122").Value = 187506432
$ws.Range("L122").Value = 200008800
$ws.Range("M122").Value = -187503982
$ws.Range("N122").Value = -200013700
$ws.Range("H126").Value = 4773.25
$ws.Range("I126").Value = 5307.7
$ws.Range("J126").Value = 3882.5
$ws.Range("K126").Value = 15923.1
$ws.Range("L126").Value = 11647.5
$ws.Range("M126").Value = -13453.1
$ws.Range("N126").Value = -16587.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 373334
$ws.Range("I81").Value = 550000
$ws.Range("J81").Value = 20002
$ws.Range("K81").Value = 1100000
$ws.Range("L81").Value = 40004
$ws.Range("M81").Value = -1098939
$ws.Range("N81").Value = -42126
$ws.Range("H84").Value = 373334
$ws.Range("I84").Value = 550000
$ws.Range("J84").Value = 20002
$ws.Range("K84").Value = 5500000
$ws.Range("L84").Value = 200020
$ws.Range("M84").Value = -5494696
$ws.Range("N84").Value = -210628
